# Reassign the weekly Fecha/Volumen/Precio/Origen figures per row (rows 2-8, 10-14
# are re-shuffled to new values; row 9 stays the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44453; J = 50; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
    3  = @{ D = 44376; J = 25; K = 18000; L = 18000; M = 18000; O = "Provincia de Limarí"; P = 600 }
    4  = @{ D = 44418; J = 30; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 500 }
    5  = @{ D = 44474; J = 45; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí"; P = 333 }
    6  = @{ D = 44446; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 467 }
    7  = @{ D = 44460; J = 45; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 433 }
    8  = @{ D = 44449; J = 45; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
    10 = @{ D = 44421; J = 25; K = 15000; L = 16000; M = 15400; O = "Provincia de Limarí"; P = 513 }
    11 = @{ D = 44432; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 467 }
    12 = @{ D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 467 }
    13 = @{ D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 467 }
    14 = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio $/Kg
}
